$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the rest of row 1 (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..32 for columns I and J
$data = @(
    @(1, 4),
    @(6, 7),
    @(1, 3),
    @(8, 8),
    @(1, 6),
    @(5, 6),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 8),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 7),
    @(1, 8),
    @(1, 5),
    @(1, 6),
    @(1, 3),
    @(1, 8),
    @(1, 1),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
